$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily OHLC + volume rows for September 2025 (end-of-month update)
$dates = @(45901, 45902, 45903, 45904, 45905, 45908, 45909, 45910, 45911, 45912, 45915, 45916, 45917, 45918, 45919, 45922, 45923, 45924, 45925, 45926, 45929, 45930)
$openV = @(4568.5341796875, 4594.89306640625, 4524.47021484375, 4543.083984375, 4363.6630859375, 4576.3798828125, 4550.2529296875, 4502.4111328125, 4527.298828125, 4669.537109375, 4702.67919921875, 4722.7548828125, 4761.9638671875, 4839.5517578125, 4826.06298828125, 4823.7568359375, 4896.8681640625, 4842.8271484375, 4966.55908203125, 5026.1640625, 4961.94091796875, 5048.3330078125)
$highV = @(4604.673828125, 4622.037109375, 4561.39306640625, 4574.67919921875, 4562.26318359375, 4606.60498046875, 4570.39501953125, 4549.59521484375, 4685.68310546875, 4713.10888671875, 4784.35498046875, 4762.72900390625, 4862.23486328125, 4927.35205078125, 4867.9658203125, 4859.26318359375, 4931.203125, 4979.001953125, 5086.27783203125, 5060.22607421875, 5072.408203125, 5108.31201171875)
$lowV  = @(4533.5517578125, 4464.365234375, 4467.47412109375, 4299.494140625, 4349.93212890625, 4521.751953125, 4470.61083984375, 4474.02490234375, 4495.130859375, 4646.84716796875, 4693.48388671875, 4675.52783203125, 4745.76708984375, 4761.69091796875, 4787.8779296875, 4785.330078125, 4769.18115234375, 4830.4072265625, 4953.06298828125, 4935.95703125, 4959.8330078125, 5034.41015625)
$closeV = @(4603.94287109375, 4505.26513671875, 4515.4091796875, 4355.9638671875, 4556.74609375, 4585.6767578125, 4494.52587890625, 4521.88818359375, 4685.68310546875, 4647.26806640625, 4720.708984375, 4756.94921875, 4853.43310546875, 4831.673828125, 4810.31396484375, 4859.26318359375, 4867.81103515625, 4979.001953125, 5047.5322265625, 4937.31201171875, 5048.05517578125, 5058.10009765625)
$volV = @(58651324, 61308764, 50863828, 55264000, 49084180, 50966412, 41471416, 39504980, 50482748, 49040820, 45102512, 45348600, 48392916, 67320032, 44922664, 39246152, 47641908, 47361232, 48754896, 43246920, 49019380, 42040372)

$startRow = 1349
$endRow = $startRow + $dates.Count - 1

for ($i = 0; $i -lt $dates.Count; $i++) {
    $r = $startRow + $i
    $prev = $r - 1
    # Clone formatting (incl. the date style on column A) from the row above
    $ws.Range("A" + $prev + ":H" + $prev).Copy($ws.Range("A" + $r + ":H" + $r))

    $ws.Range("A$r").Value = $dates[$i]
    $ws.Range("B$r").Value = $openV[$i]
    $ws.Range("C$r").Value = $highV[$i]
    $ws.Range("D$r").Value = $lowV[$i]
    $ws.Range("E$r").Value = $closeV[$i]
    $ws.Range("G$r").Value = $volV[$i]
    $ws.Range("H$r").Value = $volV[$i]
}

# Column F: fill down the E/1000 formula as one shared formula block
$ws.Range("F" + $startRow + ":F" + $endRow).Formula = "=E" + $startRow + "/1000"

# Sheet view: keep the tail of the data visible, selection parked on the next empty row
$ws.Application.ActiveWindow.ScrollRow = $startRow
$ws.Range("A" + ($endRow + 1)).Select()
